$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: B5 changes from text "101" to a real number 101
$ws.Range("B5").Value = 101

# New row 6
$ws.Range("A6").Value = "Prakash"
$ws.Range("B6").Value = 101
$ws.Range("C6").Value = "18:01:04"

# New row 7
$ws.Range("A7").Value = "Prakash"
$ws.Range("B7").Value = 101
$ws.Range("C7").Value = "18:14:35"

# New row 8
$ws.Range("A8").Value = "Prakash"
$ws.Range("B8").Value = 101
$ws.Range("C8").Value = "18:18:17"

# New row 9 - Roll Number stays a text value "103" (not numeric), like original row 5 used to be
$ws.Range("A9").Value = "Kolass"
$ws.Range("B9").Value = "'103"
$ws.Range("C9").Value = "18:18:46"
